$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.169.72'
$ws.Range('E2').Value = '  -1.83%  '

$ws.Range('D3').Value = '1.563.30'
$ws.Range('E3').Value = '  -1.56%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '206.82'
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('E6').Value = '  -1.61%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '22.34'
$ws.Range('E8').Value = '  +0.54%  '

$ws.Range('E9').Value = '  -2.01%  '

$ws.Range('E10').Value = '  +0.21%  '

$ws.Range('E11').Value = '  -0.68%  '

$ws.Range('D12').Value = '1.785.61'
$ws.Range('E12').Value = '  -1.57%  '

$ws.Range('D13').Value = '1.562.91'
$ws.Range('E13').Value = '  -1.68%  '

$ws.Range('E14').Value = '  -2.01%  '

$ws.Range('E15').Value = '  -2.29%  '

$ws.Range('D16').Value = '62.88'
$ws.Range('E16').Value = '  -0.91%  '

$ws.Range('D17').Value = '27.155.51'
$ws.Range('E17').Value = '  -1.79%  '

$ws.Range('D18').Value = '213.06'
$ws.Range('E18').Value = '  -2.76%  '

$ws.Range('E19').Value = '  -1.30%  '

$ws.Range('D20').Value = '7.22'
$ws.Range('E20').Value = '  -1.36%  '

$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('E22').Value = '  -0.14%  '

$ws.Range('E23').Value = '  -1.93%  '

$ws.Range('E24').Value = '  +0.22%  '

$ws.Range('D25').Value = '152.18'
$ws.Range('E25').Value = '  -1.01%  '

$ws.Range('D26').Value = '6.60'
$ws.Range('E26').Value = '  -3.55%  '

$ws.Range('D27').Value = '14.89'
$ws.Range('E27').Value = '  -1.45%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.104'
$ws.Range('E28').Value = '  -1.34%  '

$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.06%  '

$ws.Range('E30').Value = '  -0.74%  '

$ws.Range('E31').Value = '  -0.92%  '

$ws.Range('E32').Value = '  -1.62%  '

$ws.Range('D33').Value = '1.384.32'
$ws.Range('E33').Value = '  +0.94%  '

$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  +0.65%  '

$ws.Range('E35').Value = '  +0.52%  '

$ws.Range('D36').Value = '0.947'
$ws.Range('E36').Value = '  -2.97%  '

$ws.Range('E37').Value = '  -1.03%  '

$ws.Range('D38').Value = '0.0166'
$ws.Range('E38').Value = '  -0.96%  '

$ws.Range('D39').Value = '0.816'
$ws.Range('E39').Value = '  -0.94%  '

$ws.Range('D40').Value = '0.520'
$ws.Range('E40').Value = '  -2.84%  '

$ws.Range('E41').Value = '  -0.04%  '

$ws.Range('D42').Value = '0.991'
$ws.Range('E42').Value = '  +1.91%  '

$ws.Range('E43').Value = '  +4.33%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '63.47'
$ws.Range('E44').Value = '  -0.96%  '

$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('E46').Value = '  +0.55%  '

$ws.Range('D47').Value = '1.698.05'
$ws.Range('E47').Value = '  -1.56%  '

$ws.Range('D48').Value = '85.70'
$ws.Range('E48').Value = '  -1.91%  '

$ws.Range('D49').Value = '0.0₇0992'
$ws.Range('E49').Value = '  -1.18%  '

$ws.Range('E50').Value = '  -0.44%  '

$ws.Range('E51').Value = '  +0.19%  '
